$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'309.20"
$ws.Range("E2").Value = "'0.20%"
$ws.Range("D3").Value = "'41.33"
$ws.Range("E3").Value = "'0.83%"
$ws.Range("D4").Value = "'5.218"
$ws.Range("E4").Value = "'2.07%"
$ws.Range("E5").Value = "'0.83%"
$ws.Range("D6").Value = "'1.639"
$ws.Range("E6").Value = "'2.04%"
$ws.Range("D7").Value = "'0.9146"
$ws.Range("E7").Value = "'1.37%"
$ws.Range("D8").Value = "'2.423"
$ws.Range("E8").Value = "'-1.50%"
$ws.Range("D9").Value = "'0.1236"
$ws.Range("E9").Value = "'10.30%"
$ws.Range("D10").Value = "'0.1825"
$ws.Range("E10").Value = "'1.87%"
$ws.Range("D11").Value = "'0.09217"
$ws.Range("E11").Value = "'0.21%"
$ws.Range("D12").Value = "'0.04182"
$ws.Range("E12").Value = "'0.62%"
$ws.Range("E13").Value = "'-0.01%"
$ws.Range("D14").Value = "'0.001262"
$ws.Range("E14").Value = "'0.20%"
$ws.Range("D15").Value = "'0.005896"
$ws.Range("E15").Value = "'0.87%"
$ws.Range("E16").Value = "'0.22%"
$ws.Range("D17").Value = "'4.318"
$ws.Range("E17").Value = "'1.64%"
$ws.Range("E18").Value = "'-0.10%"
$ws.Range("D19").Value = "'7.384"
$ws.Range("E19").Value = "'11.42%"
$ws.Range("D20").Value = "'0.1402"
$ws.Range("E20").Value = "'2.71%"
$ws.Range("D21").Value = "'0.2821"
$ws.Range("E21").Value = "'4.33%"
$ws.Range("E22").Value = "'-1.33%"
$ws.Range("E23").Value = "'1.73%"
$ws.Range("D24").Value = "'0.004107"
$ws.Range("E24").Value = "'0.41%"
$ws.Range("D25").Value = "'0.0001301"
$ws.Range("E25").Value = "'0.00%"
$ws.Range("E38").Value = "'5.97%"
$ws.Range("D39").Value = "'0.05340"
$ws.Range("E39").Value = "'2.28%"
$ws.Range("D40").Value = "'0.007840"
$ws.Range("E40").Value = "'0.56%"
$ws.Range("D41").Value = "'0.1313"
$ws.Range("E41").Value = "'1.06%"
$ws.Range("D42").Value = "'0.006668"
$ws.Range("E42").Value = "'-5.45%"
$ws.Range("D43").Value = "'0.001861"
$ws.Range("E43").Value = "'-4.61%"
$ws.Range("D44").Value = "'0.007401"
$ws.Range("E44").Value = "'-3.94%"
$ws.Range("D45").Value = "'0.3066"
$ws.Range("E45").Value = "'-8.30%"
$ws.Range("D46").Value = "'0.00006763"
$ws.Range("E46").Value = "'-2.57%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'0.00%"
$ws.Range("D48").Value = "'0.2295"
$ws.Range("E48").Value = "'380.57%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'0.00%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'0.00%"
